$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6277.8696
$ws.Range("I17").Value = 1300
$ws.Range("J17").Value = 6504.136
$ws.Range("K17").Value = 3900
$ws.Range("L17").Value = 19512.408
$ws.Range("M17").Value = -3732
$ws.Range("N17").Value = -19848.408

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9249.75
$ws.Range("I62").Value = 999
$ws.Range("J62").Value = 10428.429
$ws.Range("K62").Value = 999
$ws.Range("L62").Value = 10428.429
$ws.Range("M62").Value = -375
$ws.Range("N62").Value = -11676.429

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 9249.75
$ws.Range("I65").Value = 999
$ws.Range("J65").Value = 10428.429
$ws.Range("K65").Value = 4995
$ws.Range("L65").Value = 52142.145
$ws.Range("M65").Value = -1875
$ws.Range("N65").Value = -58382.145

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2021.525
$ws.Range("I137").Value = 1508.7727
$ws.Range("J137").Value = 2648.2222
$ws.Range("K137").Value = 4526.3181
$ws.Range("L137").Value = 7944.6666
$ws.Range("M137").Value = -1976.3181
$ws.Range("N137").Value = -13044.6666

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3327515
$ws.Range("I2").Value = 4350998.5
$ws.Range("J2").Value = 1193.75
$ws.Range("K2").Value = 4350998.5
$ws.Range("L2").Value = 1193.75
$ws.Range("M2").Value = -4350885.5
$ws.Range("N2").Value = -1419.75

# ARM row 3
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 7633
$ws.Range("I3").Value = 900
$ws.Range("J3").Value = 10999.5
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 10999.5
$ws.Range("M3").Value = -785
$ws.Range("N3").Value = -11229.5

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5067.953
$ws.Range("I32").Value = 3524.7727
$ws.Range("J32").Value = 8462.950000000001
$ws.Range("K32").Value = 3524.7727
$ws.Range("L32").Value = 8462.950000000001
$ws.Range("M32").Value = -3237.7727
$ws.Range("N32").Value = -9036.950000000001

# ARM row 34
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 10000
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -10542

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 11992192
$ws.Range("I45").Value = 20552474
$ws.Range("J45").Value = 7797.6
$ws.Range("K45").Value = 20552474
$ws.Range("L45").Value = 7797.6
$ws.Range("M45").Value = -20552097
$ws.Range("N45").Value = -8551.6

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 156909.78
$ws.Range("I74").Value = 81036.62
$ws.Range("J74").Value = 354180
$ws.Range("K74").Value = 81036.62
$ws.Range("L74").Value = 354180
$ws.Range("M74").Value = -80162.62

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 156909.78
$ws.Range("I77").Value = 81036.62
$ws.Range("J77").Value = 354180
$ws.Range("K77").Value = 405183.1
$ws.Range("L77").Value = 1770900
$ws.Range("M77").Value = -400815.1

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2781434.5
$ws.Range("I102").Value = 3626525
$ws.Range("J102").Value = 4709.143
$ws.Range("K102").Value = 3626525
$ws.Range("L102").Value = 4709.143
$ws.Range("M102").Value = -3624903

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1158278.1
$ws.Range("I110").Value = 1634976.6
$ws.Range("J110").Value = 581.8570999999999
$ws.Range("K110").Value = 1634976.6
$ws.Range("L110").Value = 581.8570999999999
$ws.Range("M110").Value = -1632931.6
$ws.Range("N110").Value = -4671.8571

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3327515
$ws.Range("I116").Value = 4350998.5
$ws.Range("J116").Value = 1193.75
$ws.Range("K116").Value = 4350998.5
$ws.Range("L116").Value = 1193.75
$ws.Range("M116").Value = -4348704.5
$ws.Range("N116").Value = -5781.75

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 550820
$ws.Range("I122").Value = 2073.5715
$ws.Range("J122").Value = 2087309.9
$ws.Range("K122").Value = 6220.7145
$ws.Range("L122").Value = 6261929.699999999
$ws.Range("M122").Value = -3770.7145
$ws.Range("N122").Value = -6266829.699999999

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1991.963
$ws.Range("I132").Value = 1207.5714
$ws.Range("J132").Value = 4737.3335
$ws.Range("K132").Value = 3622.7142
$ws.Range("L132").Value = 14212.0005
$ws.Range("M132").Value = -1092.7142

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 96912.57000000001
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 96912.57000000001
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 96912.57000000001
$ws.Range("N139").Value = -107192.57

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3327515
$ws.Range("I3").Value = 4350998.5
$ws.Range("J3").Value = 1193.75
$ws.Range("K3").Value = 4350998.5
$ws.Range("L3").Value = 1193.75
$ws.Range("M3").Value = -4350884.5
$ws.Range("N3").Value = -1421.75

# BSM row 7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 501.5
$ws.Range("I7").Value = 501.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 501.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -388.5

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1885.85
$ws.Range("I20").Value = 2102.8
$ws.Range("J20").Value = 1235
$ws.Range("K20").Value = 2102.8
$ws.Range("L20").Value = 1235
$ws.Range("M20").Value = -1855.8
$ws.Range("N20").Value = -1729

# BSM row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 17099.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 17099.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 17099.5
$ws.Range("N81").Value = -19221.5

# BSM row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 17099.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 17099.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 51298.5
$ws.Range("N84").Value = -61906.5

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3908260.5
$ws.Range("I105").Value = 4809436
$ws.Range("J105").Value = 3166.6667
$ws.Range("K105").Value = 4809436
$ws.Range("L105").Value = 3166.6667
$ws.Range("M105").Value = -4807689
$ws.Range("N105").Value = -6660.6667

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2098.3
$ws.Range("I134").Value = 616.53656
$ws.Range("J134").Value = 8848.556
$ws.Range("K134").Value = 1849.60968
$ws.Range("L134").Value = 26545.668
$ws.Range("M134").Value = 685.39032

# CRP row 12
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 988.3333
$ws.Range("I12").Value = 1032.5
$ws.Range("J12").Value = 900
$ws.Range("K12").Value = 1032.5
$ws.Range("L12").Value = 900
$ws.Range("M12").Value = -862.5

# CRP row 100
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2584.913
$ws.Range("I122").Value = 2290.6667
$ws.Range("J122").Value = 3136.625
$ws.Range("K122").Value = 6872.000100000001
$ws.Range("L122").Value = 9409.875
$ws.Range("M122").Value = -4422.000100000001
$ws.Range("N122").Value = -14309.875

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 38537.56
$ws.Range("I134").Value = 47004
$ws.Range("J134").Value = 4671.8
$ws.Range("K134").Value = 141012
$ws.Range("L134").Value = 14015.4
$ws.Range("M134").Value = -138477

# CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 106
$ws.Range("I33").Value = 87.61539
$ws.Range("J33").Value = 145.83333
$ws.Range("K33").Value = 525.6923400000001
$ws.Range("L33").Value = 874.9999799999999
$ws.Range("M33").Value = -242.6923400000001
$ws.Range("N33").Value = -1440.99998

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 371518.12
$ws.Range("I46").Value = 1667058.5
$ws.Range("J46").Value = 1363.7142
$ws.Range("K46").Value = 5001175.5
$ws.Range("L46").Value = 4091.1426
$ws.Range("M46").Value = -5001084.5
$ws.Range("N46").Value = -4273.142599999999

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 74279
$ws.Range("I55").Value = 304
$ws.Range("J55").Value = 94454
$ws.Range("K55").Value = 912
$ws.Range("L55").Value = 283362
$ws.Range("M55").Value = -735
$ws.Range("N55").Value = -283716

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2903.2727
$ws.Range("I113").Value = 5275.2
$ws.Range("J113").Value = 1872
$ws.Range("K113").Value = 15825.6
$ws.Range("L113").Value = 5616
$ws.Range("M113").Value = -13655.6

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 992916.3
$ws.Range("I97").Value = 1082991.1
$ws.Range("J97").Value = 2094
$ws.Range("K97").Value = 1082991.1
$ws.Range("L97").Value = 2094
$ws.Range("M97").Value = -1082495.1
$ws.Range("N97").Value = -3086

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 15154202
$ws.Range("I113").Value = 41667696
$ws.Range("J113").Value = 3634.5715
$ws.Range("K113").Value = 41667696
$ws.Range("L113").Value = 3634.5715
$ws.Range("M113").Value = -41665526

# GSM row 119
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 90000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 90000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 90000
$ws.Range("N119").Value = -99676

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3711.6177
$ws.Range("I122").Value = 2959.08
$ws.Range("J122").Value = 5802
$ws.Range("K122").Value = 8877.24
$ws.Range("L122").Value = 17406
$ws.Range("M122").Value = -6427.24
$ws.Range("N122").Value = -22306

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8588970
$ws.Range("I126").Value = 5684619
$ws.Range("J126").Value = 11908229
$ws.Range("K126").Value = 17053857
$ws.Range("L126").Value = 35724687
$ws.Range("M126").Value = -17051387
$ws.Range("N126").Value = -35729627

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2909.8809
$ws.Range("I132").Value = 2297.4138
$ws.Range("J132").Value = 4276.154
$ws.Range("K132").Value = 6892.241399999999
$ws.Range("L132").Value = 12828.462
$ws.Range("M132").Value = -4362.241399999999

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4660
$ws.Range("I7").Value = 2790.3
$ws.Range("J7").Value = 8399.4
$ws.Range("K7").Value = 2790.3
$ws.Range("L7").Value = 8399.4
$ws.Range("M7").Value = -2678.3

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7900.467
$ws.Range("I40").Value = 4950.9
$ws.Range("J40").Value = 13799.6
$ws.Range("K40").Value = 4950.9
$ws.Range("L40").Value = 13799.6
$ws.Range("M40").Value = -4814.9

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15876143
$ws.Range("I61").Value = 15876143
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 15876143
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -15875941
$ws.Range("N61").ClearContents()

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 15876143
$ws.Range("I113").Value = 15876143
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 15876143
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -15873973
$ws.Range("N113").ClearContents()

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7092.4
$ws.Range("I122").Value = 4849.6665
$ws.Range("J122").Value = 8053.5713
$ws.Range("K122").Value = 14548.9995
$ws.Range("L122").Value = 24160.7139
$ws.Range("M122").Value = -12098.9995
$ws.Range("N122").Value = -29060.7139

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4660
$ws.Range("I126").Value = 2790.3
$ws.Range("J126").Value = 8399.4
$ws.Range("K126").Value = 8370.900000000001
$ws.Range("L126").Value = 25198.2
$ws.Range("M126").Value = -5900.900000000001

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 843.5789
$ws.Range("I113").Value = 651.53845
$ws.Range("J113").Value = 1259.6666
$ws.Range("K113").Value = 1954.61535
$ws.Range("L113").Value = 3778.9998
$ws.Range("M113").Value = 215.38465
